$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 204 (the "恐れの中に生きる者は決して自由になれない" post).
# This shifts all subsequent rows (205-212) up by one, becoming 204-211.
$ws.Rows.Item(204).Delete()
